$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the "data" column (A2:A4) from inline-string text dates into real
# Excel date serial numbers, formatted with a date-time number format.
$ws.Range("A2").Value = 45658
$ws.Range("A3").Value = 45659
$ws.Range("A4").Value = 45660

$ws.Range("A2:A4").NumberFormat = "yyyy-mm-dd h:mm:ss"
